# Vertretungsalarm.pptx update:
#  - date placeholder (master + all layouts that have one) bumped from
#    15.08.2019 -> 18.08.2019
#  - the "Fuer dich steht zur Zeit nichts..." textbox on slide 8 is
#    un-rotated and re-positioned

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16
$newDate = "18.08.2019"

function Set-DatePlaceholderText {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes

# Every slide layout's date placeholder (layout "1_Leer" has none, the
# loop just skips it since nothing matches the placeholder type there).
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Set-DatePlaceholderText $layout.Shapes
}

# Slide 8: "Textfeld 20" loses its rotation and moves up/left while
# keeping its original size.
$slide8 = $p.Slides.Item(8)
$textShape = $slide8.Shapes.Item(3)

$textShape.Rotation = 0
$textShape.Left = 389.04094488188974
$textShape.Top = 146.25803449606298
